$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.366.74'
$ws.Range('E2').Value = '  -0.05%  '

$ws.Range('D3').Value = '2.611.21'
$ws.Range('E3').Value = '  +2.83%  '

$ws.Range('E4').Value = '  +0.26%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.45%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.92%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.602'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.41%  '

$ws.Range('E8').Value = '  +0.10%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.580'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.79%  '

$ws.Range('E10').Value = '  +1.43%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0844'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.04%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.13'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.34%  '

$ws.Range('E13').Value = '  +2.27%  '

$ws.Range('D14').Value = '3.009.64'
$ws.Range('E14').Value = '  +2.87%  '

$ws.Range('E15').Value = '  +0.84%  '

$ws.Range('D16').Value = '2.614.64'
$ws.Range('E16').Value = '  +2.59%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.920'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.71%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.95'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.12%  '

$ws.Range('D19').Value = '46.485.24'
$ws.Range('E19').Value = '  +0.24%  '

$ws.Range('E20').Value = '  +1.18%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.92'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.54%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.73'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.68%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.43'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.51%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '272.75'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.09%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.05'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.40%  '

$ws.Range('E26').Value = '  +1.87%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '29.16'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +21.06%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.18%  '

$ws.Range('E29').Value = '  -0.84%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.60'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.75%  '

$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '38.77'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.64%  '

$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.21'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.05%  '

$ws.Range('E33').Value = '  +5.64%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.66'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.04%  '

$ws.Range('E35').Value = '  -1.88%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.24'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.49%  '

$ws.Range('E37').Value = '  -0.77%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '151.01'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.67%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.122'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.91%  '

$ws.Range('E40').Value = '  +1.84%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '23.41'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +34.50%  '

$ws.Range('E42').Value = '  -3.84%  '

$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0331'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.67%  '

$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.63'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.09%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.07'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.73%  '

$ws.Range('D46').Value = '2.113.56'
$ws.Range('E46').Value = '  +4.86%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.997'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.12%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '93.73'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.15%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.56'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.02%  '

$ws.Range('E50').Value = '  -4.72%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '108.98'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.66%  '
